# StaticData 외래키 테스트 (#73)
# Adds a new "MultiForeignTest" worksheet (after the existing "GroupTest"
# sheet) that exercises a multi-column / nested foreign-key lookup.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet (GroupTest) so it lands
# at the end of the tab strip, then rename + activate it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "MultiForeignTest"
[void]$ws.Activate()

# A1 documents the anchor cell of the table below (mirrors the pattern
# used on the other *Test sheets, e.g. GroupTest's A1 = "C2").
$ws.Range("A1").Value = "D5"

# Header row for the nested-foreign-key table.
$ws.Range("D5").Value = "Id"
$ws.Range("E5").Value = "TargetId"
$ws.Range("F5").Value = "Info"

# Data rows.
$ws.Range("D6").Value = 5000
$ws.Range("E6").Value = 1001
$ws.Range("F6").Value = "중첩 외래키"

$ws.Range("D7").Value = 5001
$ws.Range("E7").Value = 1003
$ws.Range("F7").Value = "테스트"

# Match the author's saved selection/view state for the new active sheet.
[void]$ws.Range("A1:G8").Select()
